# Auto-generated edit script: refreshes market-price derived columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4259.3335
$ws.Range("I74").Value = 4325.5
$ws.Range("J74").Value = 4127
$ws.Range("K74").Value = 4325.5
$ws.Range("L74").Value = 4127
$ws.Range("M74").Value = -3389.5
$ws.Range("N74").Value = -5999

$ws.Range("H77").Value = 4259.3335
$ws.Range("I77").Value = 4325.5
$ws.Range("J77").Value = 4127
$ws.Range("K77").Value = 21627.5
$ws.Range("L77").Value = 20635
$ws.Range("M77").Value = -16947.5
$ws.Range("N77").Value = -29995

$ws.Range("H100").Value = 1836
$ws.Range("I100").Value = 1393.25
$ws.Range("J100").Value = 2500.125
$ws.Range("K100").Value = 1393.25
$ws.Range("L100").Value = 2500.125
$ws.Range("M100").Value = -852.25
$ws.Range("N100").Value = -3582.125

$ws.Range("H111").Value = 1984.8
$ws.Range("I111").Value = 750
$ws.Range("J111").Value = 2433.818
$ws.Range("K111").Value = 2250
$ws.Range("L111").Value = 7301.454000000001
$ws.Range("M111").Value = 817
$ws.Range("N111").Value = -13435.454

$ws.Range("H112").Value = 41668676
$ws.Range("I112").Value = 250000580
$ws.Range("J112").Value = 2295.95
$ws.Range("K112").Value = 750001740
$ws.Range("L112").Value = 6887.849999999999
$ws.Range("M112").Value = -750000632
$ws.Range("N112").Value = -9103.849999999999

$ws.Range("H137").Value = 2603.56
$ws.Range("I137").Value = 1976
$ws.Range("J137").Value = 2956.5625
$ws.Range("K137").Value = 5928
$ws.Range("L137").Value = 8869.6875
$ws.Range("M137").Value = -3378
$ws.Range("N137").Value = -13969.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 150.4
$ws.Range("I5").Value = 162.5
$ws.Range("J5").Value = 102
$ws.Range("K5").Value = 162.5
$ws.Range("L5").Value = 102
$ws.Range("M5").Value = -50.5
$ws.Range("N5").Value = -326

$ws.Range("H32").Value = 23055
$ws.Range("I32").Value = 22636.346
$ws.Range("J32").Value = 24334.223
$ws.Range("K32").Value = 22636.346
$ws.Range("L32").Value = 24334.223
$ws.Range("M32").Value = -22349.346
$ws.Range("N32").Value = -24908.223

$ws.Range("H61").Value = 78891.69500000001
$ws.Range("I61").Value = 49086.855
$ws.Range("K61").Value = 49086.855
$ws.Range("M61").Value = -48874.855

$ws.Range("H136").Value = 78891.69500000001
$ws.Range("I136").Value = 49086.855
$ws.Range("K136").Value = 147260.565
$ws.Range("M136").Value = -144710.565

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 150.4
$ws.Range("I4").Value = 162.5
$ws.Range("J4").Value = 102
$ws.Range("K4").Value = 162.5
$ws.Range("L4").Value = 102
$ws.Range("M4").Value = -47.5
$ws.Range("N4").Value = -332

$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H20").Value = 1215.6897
$ws.Range("I20").Value = 998.8946999999999
$ws.Range("K20").Value = 998.8946999999999
$ws.Range("M20").Value = -751.8946999999999

$ws.Range("H86").Value = 14069.55
$ws.Range("I86").Value = 15249.5
$ws.Range("K86").Value = 15249.5
$ws.Range("M86").Value = -14126.5

$ws.Range("H89").Value = 14069.55
$ws.Range("I89").Value = 15249.5
$ws.Range("K89").Value = 76247.5
$ws.Range("M89").Value = -70631.5

$ws.Range("H94").Value = 1270.2
$ws.Range("I94").Value = 1350.25
$ws.Range("J94").Value = 950
$ws.Range("K94").Value = 1350.25
$ws.Range("L94").Value = 950
$ws.Range("M94").Value = -899.25
$ws.Range("N94").Value = -1852

$ws.Range("H99").Value = 1036.5
$ws.Range("I99").Value = 944.44446
$ws.Range("J99").Value = 1202.2
$ws.Range("K99").Value = 944.44446
$ws.Range("L99").Value = 1202.2
$ws.Range("M99").Value = 553.55554
$ws.Range("N99").Value = -4198.2

$ws.Range("H134").Value = 2882.6
$ws.Range("I134").Value = 2918.7917
$ws.Range("K134").Value = 8756.375100000001
$ws.Range("M134").Value = -6221.375100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 6944769
$ws.Range("I22").Value = 20833332
$ws.Range("J22").Value = 487.5
$ws.Range("K22").Value = 20833332
$ws.Range("L22").Value = 487.5
$ws.Range("M22").Value = -20832982
$ws.Range("N22").Value = -1187.5

$ws.Range("H58").Value = 2317.611
$ws.Range("I58").Value = 2351.9285
$ws.Range("J58").Value = 2197.5
$ws.Range("K58").Value = 2351.9285
$ws.Range("L58").Value = 2197.5
$ws.Range("M58").Value = -2148.9285
$ws.Range("N58").Value = -2603.5

$ws.Range("H62").Value = 3717
$ws.Range("I62").Value = 3528.3333
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3528.3333
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2904.3333
$ws.Range("N62").Value = -5248

$ws.Range("H65").Value = 3717
$ws.Range("I65").Value = 3528.3333
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 17641.6665
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -14521.6665
$ws.Range("N65").Value = -26240

$ws.Range("H107").Value = 721.36365
$ws.Range("I107").Value = 619.2857
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 619.2857
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 1300.7143
$ws.Range("N107").Value = -4740

$ws.Range("H122").Value = 2338.5264
$ws.Range("I122").Value = 2130.182
$ws.Range("K122").Value = 6390.545999999999
$ws.Range("M122").Value = -3940.545999999999

$ws.Range("H132").Value = 18349.9
$ws.Range("I132").Value = 1356.9166
$ws.Range("J132").Value = 43839.375
$ws.Range("K132").Value = 4070.7498
$ws.Range("L132").Value = 131518.125
$ws.Range("M132").Value = -1540.7498
$ws.Range("N132").Value = -136578.125

$ws.Range("H136").Value = 2317.611
$ws.Range("I136").Value = 2351.9285
$ws.Range("J136").Value = 2197.5
$ws.Range("K136").Value = 7055.7855
$ws.Range("L136").Value = 6592.5
$ws.Range("M136").Value = -4505.7855
$ws.Range("N136").Value = -11692.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 3299.1667
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 3299.1667
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 9897.500100000001
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -14797.5001

$ws.Range("H124").Value = 915
$ws.Range("I124").Value = 595.8
$ws.Range("K124").Value = 1787.4
$ws.Range("M124").Value = 3122.6

$ws.Range("H125").Value = 2857.6924
$ws.Range("I125").Value = 1670
$ws.Range("J125").Value = 3600
$ws.Range("K125").Value = 5010
$ws.Range("L125").Value = 10800
$ws.Range("M125").Value = -90
$ws.Range("N125").Value = -20640

$ws.Range("H129").Value = 2375.2593
$ws.Range("I129").Value = 1900
$ws.Range("J129").Value = 2654.8235
$ws.Range("K129").Value = 5700
$ws.Range("L129").Value = 7964.470499999999
$ws.Range("M129").Value = -700
$ws.Range("N129").Value = -17964.4705

$ws.Range("H131").Value = 1120.8182
$ws.Range("I131").Value = 422.375
$ws.Range("J131").Value = 2983.3333
$ws.Range("K131").Value = 1267.125
$ws.Range("L131").Value = 8949.999899999999
$ws.Range("M131").Value = 3772.875
$ws.Range("N131").Value = -19029.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4147.6924
$ws.Range("J80").Value = 4147.6924
$ws.Range("L80").Value = 4147.6924
$ws.Range("N80").Value = -6143.6924

$ws.Range("H83").Value = 4147.6924
$ws.Range("J83").Value = 4147.6924
$ws.Range("L83").Value = 20738.462
$ws.Range("N83").Value = -30722.462

$ws.Range("H132").Value = 75836.81
$ws.Range("I132").Value = 54148.316
$ws.Range("K132").Value = 162444.948
$ws.Range("M132").Value = -159914.948

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3968
$ws.Range("I61").Value = 3904
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3904
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3702
$ws.Range("N61").Value = -4404

$ws.Range("H68").Value = 1475
$ws.Range("I68").Value = 1416.6666
$ws.Range("J68").Value = 1562.5
$ws.Range("K68").Value = 1416.6666
$ws.Range("L68").Value = 1562.5
$ws.Range("M68").Value = -667.6666
$ws.Range("N68").Value = -3060.5

$ws.Range("H71").Value = 1475
$ws.Range("I71").Value = 1416.6666
$ws.Range("J71").Value = 1562.5
$ws.Range("K71").Value = 7083.333000000001
$ws.Range("L71").Value = 7812.5
$ws.Range("M71").Value = -3339.333000000001
$ws.Range("N71").Value = -15300.5

$ws.Range("H100").Value = 1466.963
$ws.Range("I100").Value = 1224
$ws.Range("K100").Value = 1224
$ws.Range("M100").Value = -683

$ws.Range("H113").Value = 3968
$ws.Range("I113").Value = 3904
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 3904
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -1734
$ws.Range("N113").Value = -8340

$ws.Range("H132").Value = 47927.39
$ws.Range("I132").Value = 4349.375
$ws.Range("J132").Value = 71169
$ws.Range("K132").Value = 13048.125
$ws.Range("L132").Value = 213507
$ws.Range("M132").Value = -10518.125
$ws.Range("N132").Value = -218567

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5600
$ws.Range("I62").Value = 4666.6665
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 4666.6665
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = -4042.6665
$ws.Range("N62").Value = -8248

$ws.Range("H65").Value = 5600
$ws.Range("I65").Value = 4666.6665
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 23333.3325
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = -20213.3325
$ws.Range("N65").Value = -41240

$ws.Range("H107").Value = 707.6429000000001
$ws.Range("I107").Value = 543
$ws.Range("J107").Value = 872.2857
$ws.Range("K107").Value = 1629
$ws.Range("L107").Value = 2616.8571
$ws.Range("M107").Value = 291
$ws.Range("N107").Value = -6456.8571

$ws.Range("H113").Value = 1134
$ws.Range("I113").Value = 673
$ws.Range("J113").Value = 1748.6666
$ws.Range("K113").Value = 2019
$ws.Range("L113").Value = 5245.9998
$ws.Range("M113").Value = 151
$ws.Range("N113").Value = -9585.9998

$ws.Range("H122").Value = 2431.9333
$ws.Range("I122").Value = 1589.4
$ws.Range("J122").Value = 4117
$ws.Range("K122").Value = 4768.200000000001
$ws.Range("L122").Value = 12351
$ws.Range("M122").Value = -2318.200000000001
$ws.Range("N122").Value = -17251

Write-Host "Applied Hades_Profits market data refresh (274 cell updates, 2 clears)"
